# ---------------------------------------------------------------------------
# Edit: switch the presentation's Design theme away from the custom
# "Integral" (Red Violet) theme back to the default "Office Theme" colours,
# and re-apply a (built-in) table style to the financial-documents table on
# slide 5.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# Helper: convert a "RRGGBB" hex string into the BGR-packed integer that the
# PowerPoint object model's ThemeColor.RGB / RGB() colour properties expect.
function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# --- 1. Apply the default "Office" theme colours to the deck's theme -------
# (dk1/lt1 already match; dk2, lt2 and the six accents + both hyperlink
# colours change from the "Red Violet" variant to the stock Office values.)
$officeColors = @{
    3  = "44546A"   # dk2
    4  = "E7E6E6"   # lt2
    5  = "5B9BD5"   # accent1
    6  = "ED7D31"   # accent2
    7  = "A5A5A5"   # accent3
    8  = "FFC000"   # accent4
    9  = "4472C4"   # accent5
    10 = "70AD47"   # accent6
    11 = "0563C1"   # hlink
    12 = "954F72"   # folHlink
}

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $themeColors.Item($idx).RGB = HexToComRgb $officeColors[$idx]
}

# --- 2. Re-apply the table style on the B1 financial documents table -------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{15D8F29A-6435-4967-8479-CCA74016385E}")
    }
}
